$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.871641755104065
$ws.Range("B1").Value = 2.122498035430908
$ws.Range("C1").Value = 2.311903953552246
$ws.Range("D1").Value = 3.371453046798706
$ws.Range("E1").Value = 1.432364344596863
